# Insert a new data row before existing row 229, shifting rows 229:330 down to 230:331.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("229:229").Insert()

# Populate the newly inserted row 229 with the new Cereza record.
$ws.Range("A229").Value = 9
$ws.Range("B229").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C229").Value = 'Metropolitana'
$ws.Range("D229").Value = 44875
$ws.Range("E229").Value = 13
$ws.Range("F229").Value = 'Fruta'
$ws.Range("G229").Value = 100103
$ws.Range("H229").Value = 'Frutos de hueso (carozo)'
$ws.Range("I229").Value = 100103001
$ws.Range("J229").Value = 'Cereza'
$ws.Range("K229").Value = 'Early Burlat'
$ws.Range("L229").Value = 'Primor'
$ws.Range("M229").Value = 190
$ws.Range("N229").Value = 25000
$ws.Range("O229").Value = 25000
$ws.Range("P229").Value = 25000
$ws.Range("Q229").Value = '$/bandeja 10 kilos'
$ws.Range("R229").Value = 'Provincia de Curicó'
$ws.Range("S229").Value = 2500
$ws.Range("T229").Value = 10
